# Insert a new data row at row 13 (pushing the existing rows 13-136 down to
# 14-137, each keeping its original data) and populate the new row with the
# new Mango price-report entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 13..136 down to 14..137, leaving a blank row 13 in place.
$ws.Rows("13:13").Insert()

# Fill in the new row 13 with the new record.
$ws.Range("A13").Value = 5
$ws.Range("B13").Value = "Macroferia Regional de Talca"
$ws.Range("C13").Value = "Maule"
$ws.Range("D13").Value = 44750
$ws.Range("E13").Value = 7
$ws.Range("F13").Value = "Fruta"
$ws.Range("G13").Value = 100108
$ws.Range("H13").Value = "Tropicales y subtropicales"
$ws.Range("I13").Value = 100108002
$ws.Range("J13").Value = "Mango"
$ws.Range("K13").Value = "Sin especificar"
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 200
$ws.Range("N13").Value = 8000
$ws.Range("O13").Value = 8000
$ws.Range("P13").Value = 8000
$ws.Range("Q13").Value = "$/bandeja 4 kilos"
$ws.Range("R13").Value = "Brasil"
$ws.Range("S13").Value = 2000
$ws.Range("T13").Value = 4
